$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the term/def pairs in rows 2-8 (row 1 header stays TERMS/DEF)
$ws.Range("A2").Value = "Did you hear about the claustrophobic astronaut?"
$ws.Range("B2").Value = "He just needed a little space"

$ws.Range("A3").Value = "Did you hear about the guy who stole a calendar?"
$ws.Range("B3").Value = "He got 12 months; they say his days are numbered."

$ws.Range("A4").Value = "I used to be addicted to soap,"
$ws.Range("B4").Value = " but I’m clean now."

$ws.Range("A5").Value = "I’m terrified of elevators"
$ws.Range("B5").Value = "so I’m going to start taking steps to avoid them"

$ws.Range("A6").Value = "Did you hear about the mathematician who’s afraid of negative numbers? "
$ws.Range("B6").Value = "He’ll stop at nothing to avoid them."

$ws.Range("A7").Value = "Why is that picture in jail?"
$ws.Range("B7").Value = "Because it was framed"

$ws.Range("A8").Value = "What do you call a sleeping dinosaur?"
$ws.Range("B8").Value = "Dino snore"

# New row 9 with a hyperlink in B9
$ws.Range("A9").Value = "What is a link you can't click on"
$ws.Range("B9").Value = "https://www.stylecraze.com/articles/jokes-to-tell-your-friends/#dumb-jokes-to-tell-your-friends"

[void]$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.stylecraze.com/articles/jokes-to-tell-your-friends/", "dumb-jokes-to-tell-your-friends")

[void]$ws.Range("D12").Select()
$excel.ActiveWindow.Zoom = 280
